$wb = $excel.ActiveWorkbook

# Set selection on the existing "User" sheet to F17 (was F9)
$userSheet = $wb.Worksheets.Item("User")
$userSheet.Range("F17").Select() | Out-Null

# Add a new worksheet "Library" right after "User"
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $userSheet)
$newSheet.Name = "Library"

# Populate header row and data row
$newSheet.Range("A1").Value = "name"
$newSheet.Range("B1").Value = "isbn"
$newSheet.Range("C1").Value = "author"
$newSheet.Range("A2").Value = "Java"
$newSheet.Range("B2").Value = "MKR"
$newSheet.Range("C2").Value = "Madhuri Kulkarni"

# Select D7 on the new sheet, and make it the active sheet/tab
$newSheet.Range("D7").Select() | Out-Null
